# openpyxl - manipulando as planilhas do Workbook
# Rename the existing sheet and add a brand-new blank "Sheet" after it,
# mirroring: wb.active.title = "Minha planilha"; wb.create_sheet("Sheet")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Minha planilha"

$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Sheet"
